$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.382.91"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "3.506.55"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'599.96"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").Value = "'143.14"
$ws.Range("E6").Value = "  -2.03%  "

$ws.Range("D7").Value = "3.503.45"
$ws.Range("E7").Value = "  -2.57%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = "  +2.76%  "

$ws.Range("D10").Value = "'0.133"
$ws.Range("E10").Value = "  -2.68%  "

$ws.Range("D11").Value = "'7.82"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("D12").Value = "'0.403"
$ws.Range("E12").Value = "  -3.31%  "

$ws.Range("D13").Value = "4.097.70"
$ws.Range("E13").Value = "  -2.35%  "

$ws.Range("D14").Value = "'0.0000198"
$ws.Range("E14").Value = "  -5.41%  "

$ws.Range("D15").Value = "'28.46"
$ws.Range("E15").Value = "  -5.39%  "

$ws.Range("D16").Value = "3.499.65"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "65.302.20"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").Value = "'10.98"
$ws.Range("E19").Value = "  -4.42%  "

$ws.Range("D20").Value = "'6.19"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").Value = "'14.29"
$ws.Range("E21").Value = "  -5.02%  "

$ws.Range("D22").Value = "'415.30"
$ws.Range("E22").Value = "  -4.17%  "

$ws.Range("D23").Value = "'0.597"
$ws.Range("E23").Value = "  -4.24%  "

$ws.Range("D24").Value = "'77.18"
$ws.Range("E24").Value = "  -2.59%  "

$ws.Range("D25").Value = "3.643.67"
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "'0.0000114"
$ws.Range("E27").Value = "  -5.46%  "

$ws.Range("D28").Value = "'2.43"
$ws.Range("E28").Value = "  -3.23%  "

$ws.Range("D29").Value = "'7.70"
$ws.Range("E29").Value = "  -4.90%  "

$ws.Range("D30").Value = "'8.83"
$ws.Range("E30").Value = "  -5.61%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").Value = "3.508.07"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("D33").Value = "'0.152"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("D34").Value = "'24.26"
$ws.Range("E34").Value = "  -5.07%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'7.49"
$ws.Range("E36").Value = "  -4.54%  "

$ws.Range("D37").Value = "'1.30"
$ws.Range("E37").Value = "  -10.72%  "

$ws.Range("D38").Value = "'174.77"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").Value = "'5.26"
$ws.Range("E39").Value = "  -6.89%  "

$ws.Range("D40").Value = "'1.57"
$ws.Range("E40").Value = "  -9.04%  "

$ws.Range("D41").Value = "'0.0814"
$ws.Range("E41").Value = "  -4.79%  "

$ws.Range("D42").Value = "'5.07"
$ws.Range("E42").Value = "  -3.25%  "

$ws.Range("D43").Value = "'0.852"
$ws.Range("E43").Value = "  -4.94%  "

$ws.Range("D44").Value = "'45.11"
$ws.Range("E44").Value = "  -1.95%  "

$ws.Range("D45").Value = "'1.77"
$ws.Range("E45").Value = "  -8.83%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "'2.35"
$ws.Range("E47").Value = "  -7.37%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'23.19"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.02"
$ws.Range("E49").Value = "  -2.84%  "

$ws.Range("D50").Value = "'1.08"
$ws.Range("E50").Value = "  -9.85%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "'0.900"
$ws.Range("E51").Value = "  -4.84%  "
